$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1766666666666667
$ws.Range("C2").Value = 0.55
$ws.Range("J2").Value = 0.01333333333333333
$ws.Range("P2").Value = 0.1266666666666667
$ws.Range("S2").Value = 0.1333333333333333
$ws.Range("B3").Value = 0.005847953216374269
$ws.Range("C3").Value = 0.02923976608187134
$ws.Range("J3").Value = 0.04678362573099415
$ws.Range("P3").Value = 0.7485380116959064
$ws.Range("S3").Value = 0.1695906432748538
$ws.Range("P4").Value = 0.7
$ws.Range("S4").Value = 0.3
$ws.Range("B6").Value = 0.1116071428571429
$ws.Range("F6").Value = 0.05357142857142857
$ws.Range("J6").Value = 0.2455357142857143
$ws.Range("O6").Value = 0.004464285714285714
$ws.Range("Q6").Value = 0.1919642857142857
$ws.Range("R6").Value = 0.08035714285714286
$ws.Range("S6").Value = 0.3125
$ws.Range("B7").Value = 0.1020408163265306
$ws.Range("D7").Value = 0.01020408163265306
$ws.Range("F7").Value = 0.0663265306122449
$ws.Range("J7").Value = 0.163265306122449
$ws.Range("O7").Value = 0.03061224489795918
$ws.Range("Q7").Value = 0.173469387755102
$ws.Range("R7").Value = 0.06122448979591837
$ws.Range("S7").Value = 0.3928571428571428
$ws.Range("B8").Value = 0.08126410835214447
$ws.Range("D8").Value = 0.01805869074492099
$ws.Range("E8").Value = 0.002257336343115124
$ws.Range("F8").Value = 0.05869074492099323
$ws.Range("J8").Value = 0.1309255079006772
$ws.Range("O8").Value = 0.01805869074492099
$ws.Range("Q8").Value = 0.2234762979683973
$ws.Range("R8").Value = 0.08126410835214447
$ws.Range("S8").Value = 0.3860045146726862
$ws.Range("B9").Value = 0.09375
$ws.Range("D9").Value = 0.03125
$ws.Range("E9").Value = 0.00625
$ws.Range("F9").Value = 0.04375
$ws.Range("J9").Value = 0.10625
$ws.Range("O9").Value = 0.0375
$ws.Range("Q9").Value = 0.2375
$ws.Range("R9").Value = 0.1
$ws.Range("S9").Value = 0.34375
$ws.Range("B10").Value = 0.1115326251896813
$ws.Range("D10").Value = 0.02579666160849772
$ws.Range("E10").Value = 0.0007587253414264037
$ws.Range("F10").Value = 0.07890743550834597
$ws.Range("J10").Value = 0.1206373292867982
$ws.Range("O10").Value = 0.01138088012139605
$ws.Range("Q10").Value = 0.2329286798179059
$ws.Range("R10").Value = 0.07511380880121396
$ws.Range("S10").Value = 0.3429438543247345
$ws.Range("G11").Value = 0.1275862068965517
$ws.Range("J11").Value = 0.07241379310344828
$ws.Range("K11").Value = 0.1551724137931035
$ws.Range("L11").Value = 0.6310344827586207
$ws.Range("S11").Value = 0.01379310344827586
$ws.Range("G12").Value = 0.703125
$ws.Range("J12").Value = 0.2447916666666667
$ws.Range("L12").Value = 0.04166666666666666
$ws.Range("S12").Value = 0.01041666666666667
$ws.Range("G13").Value = 0.6976744186046512
$ws.Range("J13").Value = 0.2558139534883721
$ws.Range("S13").Value = 0.04651162790697674
$ws.Range("F15").Value = 0.008658008658008658
$ws.Range("H15").Value = 0.1645021645021645
$ws.Range("I15").Value = 0.06060606060606061
$ws.Range("J15").Value = 0.3636363636363636
$ws.Range("K15").Value = 0.07792207792207792
$ws.Range("M15").Value = 0.01731601731601732
$ws.Range("O15").Value = 0.07792207792207792
$ws.Range("S15").Value = 0.2294372294372294
$ws.Range("F16").Value = 0.005181347150259068
$ws.Range("H16").Value = 0.1813471502590674
$ws.Range("I16").Value = 0.08290155440414508
$ws.Range("J16").Value = 0.3575129533678756
$ws.Range("K16").Value = 0.155440414507772
$ws.Range("M16").Value = 0.02590673575129534
$ws.Range("O16").Value = 0.05699481865284974
$ws.Range("S16").Value = 0.1347150259067358
$ws.Range("F17").Value = 0.02321083172147002
$ws.Range("H17").Value = 0.1992263056092843
$ws.Range("I17").Value = 0.06769825918762089
$ws.Range("J17").Value = 0.4468085106382979
$ws.Range("K17").Value = 0.08704061895551257
$ws.Range("M17").Value = 0.01547388781431335
$ws.Range("O17").Value = 0.05222437137330754
$ws.Range("S17").Value = 0.1083172147001934
$ws.Range("F18").Value = 0.0335195530726257
$ws.Range("H18").Value = 0.1675977653631285
$ws.Range("I18").Value = 0.09497206703910614
$ws.Range("J18").Value = 0.4245810055865922
$ws.Range("K18").Value = 0.07262569832402235
$ws.Range("M18").Value = 0.01675977653631285
$ws.Range("O18").Value = 0.0782122905027933
$ws.Range("S18").Value = 0.111731843575419
$ws.Range("F19").Value = 0.01443123938879457
$ws.Range("H19").Value = 0.2037351443123939
$ws.Range("I19").Value = 0.06706281833616298
$ws.Range("J19").Value = 0.3853989813242784
$ws.Range("K19").Value = 0.1188455008488964
$ws.Range("M19").Value = 0.02037351443123939
$ws.Range("N19").Value = 0.001697792869269949
$ws.Range("O19").Value = 0.0831918505942275
$ws.Range("S19").Value = 0.1052631578947368
